# Adds two new data rows (3 and 4) to the "Artfynd" sheet, mirroring the
# structure of the existing row 2, extending the used range to A1:AY4.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: column/row address, value kind (n=number, s=string, b=bool),
# the value itself, and whether the cell must be forced to Text format so
# that Excel does not auto-convert numeric-looking / date-looking strings.
$cells = @(
    ,@("A3", "n", 112072638, "0")
    ,@("B3", "n", 90687, "0")
    ,@("C3", "s", "Ovaliderad", "0")
    ,@("D3", "s", "LC", "0")
    ,@("E3", "n", 5964, "0")
    ,@("F3", "s", "Fjällig taggsvamp s.str.", "0")
    ,@("G3", "s", "Sarcodon imbricatus s.str.", "0")
    ,@("H3", "s", "(L.:Fr.) P.Karst.", "0")
    ,@("I3", "s", "2", "1")
    ,@("J3", "s", "fruktkroppar", "0")
    ,@("P3", "s", "Persbomossen, Upl", "0")
    ,@("Q3", "n", 654961.7553316271, "0")
    ,@("R3", "n", 6675742.431955903, "0")
    ,@("S3", "n", 4, "0")
    ,@("T3", "s", "Uppsala", "0")
    ,@("U3", "s", "Östhammar", "0")
    ,@("V3", "s", "Uppland", "0")
    ,@("W3", "s", "Dannemora", "0")
    ,@("Y3", "s", "2023-09-13", "1")
    ,@("Z3", "s", "13:42", "1")
    ,@("AA3", "s", "2023-09-13", "1")
    ,@("AB3", "s", "13:42", "1")
    ,@("AD3", "b", $false, "0")
    ,@("AE3", "b", $false, "0")
    ,@("AG3", "b", $false, "0")
    ,@("AW3", "s", "Annika Rastén", "0")
    ,@("AX3", "s", "Annika Rastén", "0")
    ,@("A4", "n", 112072636, "0")
    ,@("B4", "n", 90687, "0")
    ,@("C4", "s", "Ovaliderad", "0")
    ,@("D4", "s", "LC", "0")
    ,@("E4", "n", 5964, "0")
    ,@("F4", "s", "Fjällig taggsvamp s.str.", "0")
    ,@("G4", "s", "Sarcodon imbricatus s.str.", "0")
    ,@("H4", "s", "(L.:Fr.) P.Karst.", "0")
    ,@("I4", "s", "1", "1")
    ,@("J4", "s", "fruktkroppar", "0")
    ,@("P4", "s", "Persbomossen, Upl", "0")
    ,@("Q4", "n", 654965.1080517033, "0")
    ,@("R4", "n", 6675722.157447209, "0")
    ,@("S4", "n", 4, "0")
    ,@("T4", "s", "Uppsala", "0")
    ,@("U4", "s", "Östhammar", "0")
    ,@("V4", "s", "Uppland", "0")
    ,@("W4", "s", "Dannemora", "0")
    ,@("Y4", "s", "2023-09-13", "1")
    ,@("Z4", "s", "13:41", "1")
    ,@("AA4", "s", "2023-09-13", "1")
    ,@("AB4", "s", "13:41", "1")
    ,@("AD4", "b", $false, "0")
    ,@("AE4", "b", $false, "0")
    ,@("AG4", "b", $false, "0")
    ,@("AW4", "s", "Annika Rastén", "0")
    ,@("AX4", "s", "Annika Rastén", "0")
)

foreach ($item in $cells) {
    $addr = $item[0]
    $type = $item[1]
    $val = $item[2]
    $forceText = $item[3]
    $cell = $ws.Range($addr)
    if ($forceText -eq "1") {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $val
}

Write-Host "Rows 3 and 4 populated"
